$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Arabic rows (40-42)
$ws.Range("A40").Value = "دكني"
$ws.Range("B40").Value = "خصم 10% على جميع المنتجات"
$ws.Range("C40").Value = "A21"
$ws.Range("D40").Value = "https://dkny.com.kw/"
$ws.Range("E40").Value = "https://f.top4top.io/p_339868wzq2.png"
$ws.Range("F40").Value = "السعودية - الإمارات - الكويت"
$ws.Range("G40").Value = "لا توجد ملاحظات"

$ws.Range("A41").Value = "هواوي"
$ws.Range("B41").Value = "خصم 10% على جميع المنتجات"
$ws.Range("C41").Value = "AEB09"
$ws.Range("D41").Value = "https://consumer.huawei.com/ae-en/offer/"
$ws.Range("E41").Value = "https://g.top4top.io/p_3398epsnm3.png"
$ws.Range("F41").Value = "الإمارات"
$ws.Range("G41").Value = "لا توجد ملاحظات"

$ws.Range("A42").Value = "ريبوك"
$ws.Range("B42").Value = "خصم 15% على جميع المنتجات"
$ws.Range("C42").Value = "ADM84"
$ws.Range("D42").Value = "https://www.reebok.ae/"
$ws.Range("E42").Value = "https://e.top4top.io/p_3398jj9fc1.png"
$ws.Range("F42").Value = "الإمارات"
$ws.Range("G42").Value = "لا توجد ملاحظات"

# English rows (97-99)
$ws.Range("A97").Value = "DKNY"
$ws.Range("B97").Value = "خصم 10% على جميع المنتجات"
$ws.Range("C97").Value = "A21"
$ws.Range("D97").Value = "https://dkny.com.kw/"
$ws.Range("E97").Value = "https://f.top4top.io/p_339868wzq2.png"
$ws.Range("F97").Value = "السعودية - الإمارات - الكويت"
$ws.Range("G97").Value = "لا توجد ملاحظات"

$ws.Range("A98").Value = "Huawei"
$ws.Range("B98").Value = "خصم 10% على جميع المنتجات"
$ws.Range("C98").Value = "AEB09"
$ws.Range("D98").Value = "https://consumer.huawei.com/ae-en/offer/"
$ws.Range("E98").Value = "https://g.top4top.io/p_3398epsnm3.png"
$ws.Range("F98").Value = "الإمارات"
$ws.Range("G98").Value = "لا توجد ملاحظات"

$ws.Range("A99").Value = "Reebok"
$ws.Range("B99").Value = "خصم 15% على جميع المنتجات"
$ws.Range("C99").Value = "ADM84"
$ws.Range("D99").Value = "https://www.reebok.ae/"
$ws.Range("E99").Value = "https://e.top4top.io/p_3398jj9fc1.png"
$ws.Range("F99").Value = "الإمارات"
$ws.Range("G99").Value = "لا توجد ملاحظات"

# New empty styled rows 100-102, copying style/format from row 99's row-band (row 98 used as template of blank bordered row)
$ws.Range("A100:G102").Value = ""

# Update selection/view to match the final state
$ws.Range("F104").Select()
